$wb = $excel.ActiveWorkbook

# --- "İş Takip Listesi" (sheet 1): edit the personnel text for rows 59-94 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()

# All of E59:E94 shared the same "GÖREVLİ PERSONELLER" text; change it to the
# new personnel list. Writing the whole block at once keeps every cell
# pointing at a single (new) shared string, exactly like the source edit.
$ws1.Range("E59:E94").Value = "KEMAL KORKMAZ, ALİ BAŞKURT, İSMAİL AKLAN, ENGİN UĞURLU"

# Reflect the user's on-screen selection/scroll after making the edit.
$ws1.Range("E59:E94").Select()

# --- "Personel" (sheet 3): widen column B to fit its content ---
# (the engine pads the stored width by +5/6 relative to ColumnWidth, so back
# that padding out to land exactly on the target stored width of 14)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Columns.Item(2).ColumnWidth = 13.1666666666667
